# Add Team 6 (Database) final requirements rows to Sheet1.
# Mirrors the author's commit "added team 6 to final req xcel file":
# appends rows 107-126 (DB_1..DB_20) below the existing requirements table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@(107, "DB_1", "6 - Database", "The database shall ensure that each username is unique and associated with a UUID.", "initializeDatabase() ")
  ,@(108, "DB_2", "6 - Database", "The database shall generate a UUID for each new user record upon creation.", "initializeDatabase() ")
  ,@(109, "DB_3", "6 - Database", "The database shall provide an interface that allows the page manager to request and retrieve stored user information by username, UUID, or other attributes.", "getPlayer(String username, String password)")
  ,@(110, "DB_4", "6 - Database", "The database shall support deletion of user records manually if needed.", "deprecated")
  ,@(111, "DB_5", "6 - Database", "The database shall maintain data and value ranges for all fields to maintain consistency.", "getPlayer(String username, String password)")
  ,@(112, "DB_6", "6 - Database", "The database shall respond to requests from page manager and deliver the requested information within a short span of time.", "getPlayer(String username, String password)")
  ,@(113, "DB_7", "6 - Database", "The database shall ensure consistency by preventing duplicate names by using the UUID's.", "addPlayer(String username, String password)")
  ,@(114, "DB_8", "6 - Database", "The database shall be scalable to handle an increasing number of user records.", "addPlayer(String username, String password)")
  ,@(115, "DB_9", "6 - Database", "The database shall support the insertion of new user records.", "addPlayer(String username, String password)")
  ,@(116, "DB_10", "6 - Database", "The database shall support the update of a user's win count.", "recordMatchResult(int winnerId, int loserId)")
  ,@(117, "DB_11", "6 - Database", "The database shall support the update of a user's loss count.", "recordMatchResult(int winnerId, int loserId)")
  ,@(118, "DB_12", "6 - Database", "The database shall support the update of a user's total games played.", "recordMatchResult(int winnerId, int loserId)")
  ,@(119, "DB_13", "6 - Database", "The database shall support the update of a user's elo value.", "getPlayer(String username, String password)")
  ,@(120, "DB_14", "6 - Database", "The database shall provide a user's elo value to pair up system upon request.", "getPlayer(String username, String password)")
  ,@(121, "DB_15", "6 - Database", "The database shall store player's username.", "updatePlayerStats(int playerId, int wins, int losses, int ELO, int gamesPlayed)")
  ,@(122, "DB_16", "6 - Database", "The database shall store player's UUID.", "updatePlayerStats(int playerId, int wins, int losses, int ELO, int gamesPlayed)")
  ,@(123, "DB_17", "6 - Database", "The database shall store player's wins.", "updatePlayerStats(int playerId, int wins, int losses, int ELO, int gamesPlayed)")
  ,@(124, "DB_18", "6 - Database", "The database shall store player's losses.", "updatePlayerStats(int playerId, int wins, int losses, int ELO, int gamesPlayed)")
  ,@(125, "DB_19", "6 - Database", "The database shall store player's total games.", "updatePlayerStats(int playerId, int wins, int losses, int ELO, int gamesPlayed)")
  ,@(126, "DB_20", "6 - Database", "The database shall store player's elo.", "updatePlayerStats(int playerId, int wins, int losses, int ELO, int gamesPlayed)")
)

foreach ($row in $rows) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)

    $cellA.Value = $row[1]
    $cellB.Value = $row[2]
    $cellC.Value = $row[3]
    $cellD.Value = $row[4]

    # Columns A, B, D use the Times New Roman 12pt wrap-text style already
    # used throughout the requirements table (style index 1 in styles.xml).
    foreach ($cell in @($cellA, $cellB, $cellD)) {
        $cell.Font.Name = "Times New Roman"
        $cell.Font.Size = 12
        $cell.WrapText = $true
    }
}

# Column C (requirement text) keeps the sheet's default font with wrap text,
# matching the style already on C107 (style index 7 in styles.xml). Turn on
# wrap for the first cell, then fan that exact format out to the rest of the
# new column-C cells so we don't spawn a near-duplicate font/style entry.
$ws.Cells.Item(107, 3).WrapText = $true
$ws.Range("C107").Copy()
$ws.Range("C108:C126").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the view roughly where the author left it after entering the data.
$ws.Range("F107").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 95
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
